$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume/Hora columns to remain text so values round-trip exactly as typed
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "260.81"
$ws.Range("E2").Value = "1.70%"
$ws.Range("G2").Value = "6"

# Row 3
$ws.Range("D3").Value = "27.30"
$ws.Range("E3").Value = "1.08%"
$ws.Range("G3").Value = "6"

# Row 4
$ws.Range("D4").Value = "4.751"
$ws.Range("E4").Value = "8.79%"
$ws.Range("G4").Value = "6"

# Row 5
$ws.Range("D5").Value = "0.06068"
$ws.Range("G5").Value = "6"

# Row 6
$ws.Range("D6").Value = "6.653"
$ws.Range("E6").Value = "0.29%"
$ws.Range("G6").Value = "6"

# Row 7
$ws.Range("D7").Value = "0.8496"
$ws.Range("E7").Value = "-0.28%"
$ws.Range("G7").Value = "6"

# Row 8
$ws.Range("D8").Value = "0.9210"
$ws.Range("E8").Value = "-1.91%"
$ws.Range("G8").Value = "6"

# Row 9
$ws.Range("D9").Value = "0.1404"
$ws.Range("E9").Value = "1.51%"
$ws.Range("G9").Value = "6"

# Row 10
$ws.Range("D10").Value = "0.04894"
$ws.Range("E10").Value = "-0.30%"
$ws.Range("G10").Value = "6"

# Row 11
$ws.Range("E11").Value = "0.30%"
$ws.Range("G11").Value = "6"

# Row 12
$ws.Range("D12").Value = "0.03114"
$ws.Range("E12").Value = "1.28%"
$ws.Range("G12").Value = "6"

# Row 13
$ws.Range("D13").Value = "0.09083"
$ws.Range("E13").Value = "-0.28%"
$ws.Range("G13").Value = "6"

# Row 14
$ws.Range("D14").Value = "0.001531"
$ws.Range("E14").Value = "-0.36%"
$ws.Range("G14").Value = "6"

# Row 15
$ws.Range("D15").Value = "0.0006077"
$ws.Range("E15").Value = "-94.16%"
$ws.Range("G15").Value = "6"

# Row 16
$ws.Range("D16").Value = "0.006131"
$ws.Range("E16").Value = "-0.63%"
$ws.Range("G16").Value = "6"

# Row 17
$ws.Range("D17").Value = "3.452"
$ws.Range("E17").Value = "-0.93%"
$ws.Range("G17").Value = "6"

# Row 18
$ws.Range("D18").Value = "3.154"
$ws.Range("E18").Value = "-0.47%"
$ws.Range("G18").Value = "6"

# Row 19
$ws.Range("E19").Value = "-1.26%"
$ws.Range("G19").Value = "6"

# Row 20
$ws.Range("E20").Value = "2.49%"
$ws.Range("G20").Value = "6"

# Row 21
$ws.Range("G21").Value = "6"

# Row 22
$ws.Range("D22").Value = "4.091"
$ws.Range("E22").Value = "4.45%"
$ws.Range("G22").Value = "6"

# Row 23
$ws.Range("D23").Value = "0.04237"
$ws.Range("E23").Value = "-0.65%"
$ws.Range("G23").Value = "6"

# Row 24
$ws.Range("E24").Value = "-0.24%"
$ws.Range("G24").Value = "6"

# Row 25
$ws.Range("E25").Value = "-8.63%"
$ws.Range("G25").Value = "6"

# Row 26
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").Value = "0.06%"
$ws.Range("G26").Value = "6"

# Row 27
$ws.Range("E27").Value = "3.13%"
$ws.Range("G27").Value = "6"

# Row 28
$ws.Range("G28").Value = "6"

# Row 29
$ws.Range("G29").Value = "6"

# Row 30
$ws.Range("G30").Value = "6"

# Row 31
$ws.Range("G31").Value = "6"

# Row 32
$ws.Range("G32").Value = "6"

# Row 33
$ws.Range("G33").Value = "6"

# Row 34
$ws.Range("G34").Value = "6"

# Row 35
$ws.Range("G35").Value = "6"

# Row 36
$ws.Range("G36").Value = "6"

# Row 37
$ws.Range("G37").Value = "6"

# Row 38
$ws.Range("G38").Value = "6"

# Row 39
$ws.Range("G39").Value = "6"

# Row 40
$ws.Range("D40").Value = "0.03873"
$ws.Range("E40").Value = "1.55%"
$ws.Range("G40").Value = "6"

# Row 41
$ws.Range("D41").Value = "0.1113"
$ws.Range("E41").Value = "1.30%"
$ws.Range("G41").Value = "6"

# Row 42
$ws.Range("D42").Value = "0.004133"
$ws.Range("E42").Value = "-33.77%"
$ws.Range("G42").Value = "6"

# Row 43
$ws.Range("D43").Value = "0.01498"
$ws.Range("E43").Value = "7.19%"
$ws.Range("G43").Value = "6"

# Row 44
$ws.Range("D44").Value = "0.002209"
$ws.Range("E44").Value = "0.41%"
$ws.Range("G44").Value = "6"

# Row 45
$ws.Range("D45").Value = "0.00005321"
$ws.Range("E45").Value = "-1.04%"
$ws.Range("G45").Value = "6"

# Row 46
$ws.Range("E46").Value = "0.05%"
$ws.Range("G46").Value = "6"

# Row 47
$ws.Range("E47").Value = "-4.26%"
$ws.Range("G47").Value = "6"

# Row 48
$ws.Range("G48").Value = "6"

# Row 49
$ws.Range("E49").Value = "0.05%"
$ws.Range("G49").Value = "6"

# Row 50
$ws.Range("E50").Value = "0.05%"
$ws.Range("G50").Value = "6"

# Row 51
$ws.Range("G51").Value = "6"
